$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.682.35"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.281.60"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.47"
$ws.Range("E5").Value = "  +10.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.53"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.70"
$ws.Range("E10").Value = "  +5.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.87"
$ws.Range("E12").Value = "  +10.11%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.73"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "2.623.62"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.876"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "2.275.18"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "43.530.61"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  +12.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.09"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.41"
$ws.Range("E22").Value = "  -3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("E23").Value = "  +7.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.73"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.59"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.61"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.40"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.34"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.53"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0914"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.65"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  +5.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.92"
$ws.Range("E40").Value = "  +22.10%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.56"
$ws.Range("E41").Value = "  +14.98%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.44"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.238"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  +19.20%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.68"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.02"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0997"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.456"
$ws.Range("E51").Value = "  +2.51%  "
